{"js": "// The run \"p026r_1\" is unique in the document (inside \"<id>p026r_1</id>\",\n// itself split across three runs: \"<id>\", \"p026r_1\", \"</id>\"). Search for\n// the full visible text and replace it in place - Word merges the matched\n// range into a single run, inheriting the formatting of the first run in\n// the hit (Courier New, color 7f6000, etc.), which is exactly the target\n// of the edit (the three runs collapse into one run with the concatenated\n// text).\nconst body = context.document.body;\nconst results = body.search(\"<id>p026r_1</id>\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find '<id>p026r_1</id>' in the document body.\");\n}\n\nresults.items[0].insertText(\"<id>p026r_1</id>\", \"Replace\");\nawait context.sync();\n", "ps1": "# The visible text \"<id>p026r_1</id>\" is split across three runs in the\n# OOXML: \"<id>\", \"p026r_1\", \"</id>\". \"p026r_1\" is unique in the whole\n# document, so it is used as the anchor to locate the exact spot. The three\n# runs are merged into a single run (keeping the first run's formatting/rsid\n# attributes) by deleting the trailing \"p026r_1</id>\" text and re-appending\n# the full \"p026r_1</id>\" string onto the end of the \"<id>\" run - this is\n# equivalent to Word's own \"type over/merge adjacent runs\" behavior and\n# matches the target edit exactly (three runs -> one run with the\n# concatenated text).\n$d = $word.ActiveDocument\n\n$prefixText = \"<id>\"\n$coreText = \"p026r_1\"\n$suffixText = \"</id>\"\n\n$anchorRange = $d.Content\n$found = $anchorRange.Find.Execute($coreText)\nif (-not $found) {\n    throw \"Could not find anchor text '$coreText' in the document.\"\n}\n\n$coreStart = $anchorRange.Start\n$coreEnd = $anchorRange.End\n\n$beforeRange = $d.Range($coreStart - $prefixText.Length, $coreStart)\n$afterRange = $d.Range($coreEnd, $coreEnd + $suffixText.Length)\n\nif ($beforeRange.Text -ne $prefixText) {\n    throw \"Expected '$prefixText' immediately before '$coreText', found '$($beforeRange.Text)'.\"\n}\nif ($afterRange.Text -ne $suffixText) {\n    throw \"Expected '$suffixText' immediately after '$coreText', found '$($afterRange.Text)'.\"\n}\n\n# Remove \"p026r_1</id>\" (keeps the \"<id>\" run's range/formatting intact),\n# then grow that run with the merged text.\n$tailRange = $d.Range($coreStart, $coreEnd + $suffixText.Length)\n$tailRange.Delete()\n$beforeRange.InsertAfter($coreText + $suffixText)\n"}
